$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1215
$ws.Range("M2").Value = 5182

$ws.Range("L3").Value = 982
$ws.Range("M3").Value = 5142

$ws.Range("L4").Value = 408
$ws.Range("M4").Value = 5078

$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 5006
